# Insert a new row at row 205, pushing existing rows 205:334 down to 206:335,
# then populate the newly inserted row 205 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 205 (rows below shift down by one)
$ws.Rows.Item(205).Insert()

# Populate the new row 205 with the new data record
$ws.Cells.Item(205, 1).Value = 9
$ws.Cells.Item(205, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(205, 3).Value = "Metropolitana"
$ws.Cells.Item(205, 4).Value = 44582
$ws.Cells.Item(205, 4).Style = $ws.Cells.Item(206, 4).Style
$ws.Cells.Item(205, 4).NumberFormat = $ws.Cells.Item(206, 4).NumberFormat
$ws.Cells.Item(205, 5).Value = 13
$ws.Cells.Item(205, 6).Value = 100112039
$ws.Cells.Item(205, 7).Value = "Ciboulette"
$ws.Cells.Item(205, 8).Value = "Sin especificar"
$ws.Cells.Item(205, 9).Value = "Primera"
$ws.Cells.Item(205, 10).Value = 160
$ws.Cells.Item(205, 11).Value = 1000
$ws.Cells.Item(205, 12).Value = 1200
$ws.Cells.Item(205, 13).Value = 1100
$ws.Cells.Item(205, 14).Value = "`$/docena de atados"
$ws.Cells.Item(205, 15).Value = "Región Metropolitana"
$ws.Cells.Item(205, 16).Value = 367
$ws.Cells.Item(205, 17).Value = 3
$ws.Cells.Item(205, 18).Value = "Hortaliza"
